$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows starting at row 10 to make room for new exposure sites
$ws.Range("A10:A20").EntireRow.Insert()

# Rewrite rows 10 through 42 with the final Location/Site/Exposure period/Notes data
$ws.Range("A10").Value = "Cheltenham"
$ws.Range("B10").Value = "Bodero Southland Shopping Centre, 1239 Nepean Hwy"
$ws.Range("C10").Value = "22/12/20 6.45pm - 7pm"
$ws.Range("D10").Value = "Case shopped in store"
$ws.Range("A11").Value = "Cheltenham"
$ws.Range("B11").Value = "Chemist Warehouse Cheltenham, 326/330 Charman Rd"
$ws.Range("C11").Value = "03/01/21, 3.30pm - 3.45pm"
$ws.Range("D11").Value = "Case shopped in store"
$ws.Range("A12").Value = "Cheltenham"
$ws.Range("B12").Value = "Coles, Westfield Southland"
$ws.Range("C12").Value = "22/12/20 11:50am-12:10pm"
$ws.Range("D12").Value = "Case shopped in store"
$ws.Range("A13").Value = "Cheltenham"
$ws.Range("B13").Value = "Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Range("C13").Value = "22/12/2020 3:50pm-4:05pm"
$ws.Range("D13").Value = "Case shopped in store"
$ws.Range("A14").Value = "Cheltenham"
$ws.Range("B14").Value = "Kmart Southland Shopping Centre, 1239 Nepean Highway"
$ws.Range("C14").Value = "22/12/20 6.30pm - 6.45pm"
$ws.Range("D14").Value = "Case shopped in store"
$ws.Range("A15").Value = "Cheltenham"
$ws.Range("B15").Value = "Kmart Southland Shopping Centre, 1239 Nepean Highway"
$ws.Range("C15").Value = "28/12/20 2.30pm-3pm"
$ws.Range("D15").Value = "Case shopped in store"
$ws.Range("A16").Value = "Cheltenham"
$ws.Range("B16").Value = "Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Range("C16").Value = "22/12/2020 3:30pm-3:50pm"
$ws.Range("D16").Value = "Case shopped in store"
$ws.Range("A17").Value = "Cheltenham"
$ws.Range("B17").Value = "Specsavers, 1004-1005 Westfield Southland"
$ws.Range("C17").Value = "22/12/20 11:00am-1145am"
$ws.Range("D17").Value = "Case shopped in store"
$ws.Range("A18").Value = "Dandenong"
$ws.Range("B18").Value = "Kmart - Clayton Plaza, 2107 Dandenong Rd"
$ws.Range("C18").Value = "30/12/20 7pm - 7.30pm"
$ws.Range("D18").Value = "Case shopped at store"
$ws.Range("A19").Value = "Dandenong"
$ws.Range("B19").Value = "Woolworths - Clayton Plaza, 2107 Dandenong Rd"
$ws.Range("C19").Value = "30/12/20 7.30pm - 745pm"
$ws.Range("D19").Value = "Case shopped at store"
$ws.Range("A20").Value = "Forest Hill"
$ws.Range("B20").Value = "Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131"
$ws.Range("C20").Value = "28/12/20 12:00pm-2:00pm"
$ws.Range("D20").Value = "1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min"
$ws.Range("A21").Value = "Fountain Gate Shopping Centre"
$ws.Range("B21").Value = "Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805"
$ws.Range("C21").Value = "26/12/20 9:00am-11:00am"
$ws.Range("D21").Value = ""
$ws.Range("A22").Value = "Frankston"
$ws.Range("B22").Value = "TK Maxx Frankston, 10 Shannon Street, Bayside Shopping Centre"
$ws.Range("C22").Value = "31/12/20 2pm - 3pm"
$ws.Range("D22").Value = "Case shopped at store"
$ws.Range("A23").Value = "Glen Waverley"
$ws.Range("B23").Value = "Mocha Jos  87 Kingsway, Glen Waverley VIC 3150"
$ws.Range("C23").Value = "28/12/20 1:30pm-1:45pm"
$ws.Range("D23").Value = ""
$ws.Range("A24").Value = "Hallam"
$ws.Range("B24").Value = "Coles Hallam  2 Princes Domain Drive, Hallam VIC 3803"
$ws.Range("C24").Value = "30/12/20 6:15am-6:30am"
$ws.Range("D24").Value = "Case shopped in store"
$ws.Range("A25").Value = "Lakes Entrance"
$ws.Range("B25").Value = "Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909"
$ws.Range("C25").Value = "29/12/20 11:15am-12:15pm"
$ws.Range("D25").Value = "Case hired a boat"
$ws.Range("A26").Value = "Lakes Entrance"
$ws.Range("B26").Value = "Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909"
$ws.Range("C26").Value = "30/12/20 5:00pm-6:30pm"
$ws.Range("D26").Value = "Case attended outside premises"
$ws.Range("A27").Value = "Lakes Entrance"
$ws.Range("B27").Value = "Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909"
$ws.Range("C27").Value = "30/12/20 11:15am-11:20am"
$ws.Range("D27").Value = "Case picked up takeaway"
$ws.Range("A28").Value = "Lakes Entrance"
$ws.Range("B28").Value = "Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909"
$ws.Range("C28").Value = "30/12/20 6:00pm-6:15pm"
$ws.Range("D28").Value = "Case shopped in store"
$ws.Range("A29").Value = "Mentone"
$ws.Range("B29").Value = "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194"
$ws.Range("C29").Value = "29/12/20 07:30am-08:00am"
$ws.Range("D29").Value = "Case shopped in store"
$ws.Range("A30").Value = "Mentone"
$ws.Range("B30").Value = "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194"
$ws.Range("C30").Value = "31/12/20 08:00am-08:30am"
$ws.Range("D30").Value = "Case shopped in store"
$ws.Range("A31").Value = "Mentone"
$ws.Range("B31").Value = "Mentone/Parkdale Beach"
$ws.Range("C31").Value = "27/12/20 10:00am-4:30pm"
$ws.Range("D31").Value = ""
$ws.Range("A32").Value = "Mentone"
$ws.Range("B32").Value = "Woolworths Mentone  105-111 Balcombe Road, Mentone VIC 3194"
$ws.Range("C32").Value = "23/12/20 2:45pm-3:05pm"
$ws.Range("D32").Value = "Case shopped in store"
$ws.Range("A33").Value = "Moorabbin"
$ws.Range("B33").Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Range("C33").Value = "30/12/20 10:45am-12:15pm"
$ws.Range("D33").Value = "Case shopped in store"
$ws.Range("A34").Value = "Moorabbin"
$ws.Range("B34").Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Range("C34").Value = "30/12/20 4:00m- 5:50pm"
$ws.Range("D34").Value = "Case shopped in store"
$ws.Range("A35").Value = "Mordialloc"
$ws.Range("B35").Value = "Woodlands Golf Club  109 White Street Mordialloc VIC 3195"
$ws.Range("C35").Value = "23/12/20 8:00am-2:00pm"
$ws.Range("D35").Value = "Case attended course"
$ws.Range("A36").Value = "Mordialloc"
$ws.Range("B36").Value = "Woodlands Golf Club  109 White Street Mordialloc VIC 3195"
$ws.Range("C36").Value = "28/12/20 12:00pm-6:00pm"
$ws.Range("D36").Value = "Case attended course"
$ws.Range("A37").Value = "Mount Martha"
$ws.Range("B37").Value = "Mount Martha Fine Foods, 34 Lochiel Ave"
$ws.Range("C37").Value = "31/12/20 3pm - 3.15pm"
$ws.Range("D37").Value = "Takeaway coffee"
$ws.Range("A38").Value = "Mount Waverley"
$ws.Range("B38").Value = "Ritchies IGA  283 Stephensons Road, Mount Waverley VIC 3149"
$ws.Range("C38").Value = "30/12/20 2:00pm-2:30pm"
$ws.Range("D38").Value = "Case shopped for half an hour"
$ws.Range("A39").Value = "Oakleigh"
$ws.Range("B39").Value = "Bunnings Oakleigh  1041 Centre Road, Oakleigh South"
$ws.Range("C39").Value = "30/12/20 11:00am-11:30am"
$ws.Range("D39").Value = "Case shopped for 30 minutes"
$ws.Range("A40").Value = "Oakleigh"
$ws.Range("B40").Value = "Katialo restaurant  8 Eaton Mall, Oakleigh VIC 3166"
$ws.Range("C40").Value = "28/12/20 7:00pm-7:10pm"
$ws.Range("D40").Value = ""
$ws.Range("A41").Value = "Springvale"
$ws.Range("B41").Value = "IKEA Springvale, 917 Princes Hwy"
$ws.Range("C41").Value = "29/12/20 4pm - 6pm"
$ws.Range("D41").Value = "Case shopped at store and dined at cafe"
$ws.Range("A42").Value = "Wonthaggi"
$ws.Range("B42").Value = "Wonthaggi Plaza Shopping centre  2 Biggs Drive, Wonthaggi VIC 3995"
$ws.Range("C42").Value = "28/12/20 1:30pm-2.30pm"
$ws.Range("D42").Value = "Kmart- shopped for 15 mins"
